# Household member form: bump form_version setting and leave the
# "settings" sheet active/selected at the edited cell (matches the
# author's workflow of opening the settings sheet, editing the
# form_version value, and saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# settings!B3 is the value for the "form_version" setting (A3="form_version").
# Bump it from the placeholder 1 to the real version stamp.
$ws.Range("B3").Value = 20130408

# Make "settings" the active/selected sheet with B3 as the selected cell,
# matching the saved view state in the edited workbook.
$ws.Activate()
$ws.Range("B3").Select()
